$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update model name labels in column A
$ws.Range("A2").Value = "model_24_8_0"
$ws.Range("A3").Value = "model_24_8_22"
$ws.Range("A4").Value = "model_24_8_21"
$ws.Range("A5").Value = "model_24_8_20"
$ws.Range("A6").Value = "model_24_8_19"
$ws.Range("A7").Value = "model_24_8_18"
$ws.Range("A8").Value = "model_24_8_17"
$ws.Range("A9").Value = "model_24_8_16"
$ws.Range("A10").Value = "model_24_8_15"
$ws.Range("A11").Value = "model_24_8_14"
$ws.Range("A12").Value = "model_24_8_13"
$ws.Range("A13").Value = "model_24_8_23"
$ws.Range("A14").Value = "model_24_8_12"
$ws.Range("A15").Value = "model_24_8_10"
$ws.Range("A16").Value = "model_24_8_9"
$ws.Range("A17").Value = "model_24_8_8"
$ws.Range("A18").Value = "model_24_8_7"
$ws.Range("A19").Value = "model_24_8_6"
$ws.Range("A20").Value = "model_24_8_5"
$ws.Range("A21").Value = "model_24_8_4"
$ws.Range("A22").Value = "model_24_8_3"
$ws.Range("A23").Value = "model_24_8_2"
$ws.Range("A24").Value = "model_24_8_1"
$ws.Range("A25").Value = "model_24_8_11"
$ws.Range("A26").Value = "model_24_8_24"

# Update metrics columns B-Q with new values (same across all rows)
$ws.Range("B2:B26").Value = [double]"0.999999900753416"
$ws.Range("C2:C26").Value = [double]"0.6876245000598846"
$ws.Range("D2:D26").Value = [double]"0.9999996622279347"
$ws.Range("E2:E26").Value = [double]"0.9999994525339465"
$ws.Range("F2:F26").Value = [double]"0.9999995618376588"
$ws.Range("G2:G26").Value = [double]"5.891705945054375e-08"
$ws.Range("H2:H26").Value = [double]"0.1854395906834327"
$ws.Range("I2:I26").Value = [double]"4.960264897052699e-08"
$ws.Range("J2:J26").Value = [double]"2.203587190523217e-07"
$ws.Range("K2:K26").Value = [double]"1.349806840114243e-07"
$ws.Range("L2:L26").Value = [double]"9.931498004524596e-05"
$ws.Range("M2:M26").Value = [double]"0.0002427283655664161"
$ws.Range("N2:N26").Value = [double]"1.000000140112824"
$ws.Range("O2:O26").Value = [double]"0.0002530618146055635"
$ws.Range("P2:P26").Value = [double]"115.2942703081797"
$ws.Range("Q2:Q26").Value = [double]"165.268179127776"
